# Add a new data row (row 4) to the active worksheet, mirroring the
# existing rows' layout (header row 1, data rows 2-3 already present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 112088271
$ws.Range("B4").Value = 98446
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 222771
$ws.Range("F4").Value = "Svart trolldruva"
$ws.Range("G4").Value = "Actaea spicata"
$ws.Range("H4").Value = "L."
$ws.Range("I4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("P4").Value = "S om Idegransreservatet, Öl"
$ws.Range("Q4").Value = 619889.018100369
$ws.Range("R4").Value = 6346803.221470654
$ws.Range("S4").Value = 25
$ws.Range("T4").Value = "Kalmar"
$ws.Range("U4").Value = "Borgholm"
$ws.Range("V4").Value = "Öland"
$ws.Range("W4").Value = "Böda"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2023-09-14"
$ws.Range("Y4").Style = "Normal"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2023-09-14"
$ws.Range("AA4").Style = "Normal"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AT4").Value = ""
$ws.Range("AW4").Value = "Pontus Axén"
$ws.Range("AX4").Value = "Pontus Axén"
$ws.Range("AY4").Value = ""
